$d = $word.ActiveDocument

# The paragraph currently reads "Version 2." split across several runs,
# with an empty "_GoBack" bookmark sitting between the " 2" run and the
# trailing "." run. The target text is "Version 1." with the "Versi"/"on"
# runs merged into one "Version" run, the " 2"/"." runs merged into a
# single " 1." run, and the (empty) "_GoBack" bookmark relocated to after
# that merged run (i.e. at the very end of the paragraph's text).

# Remove the existing bookmark; we'll recreate it in the right spot once
# the text has settled.
$d.Bookmarks("_GoBack").Delete()

# Merge "Versi" + "on" into a single "Version" run.
$d.Content.Find.Execute("Version", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Version", 2) | Out-Null

# Merge " 2" + "." into a single " 1." run (this also collapses the run
# break where the old bookmark used to sit).
$d.Content.Find.Execute("2.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1.", 2) | Out-Null

$paraEnd = $d.Paragraphs(1).Range.End - 1   # end of text, before the pilcrow

# Adding a bookmark exactly at the end of the document's last paragraph is
# mishandled by this host (it snaps back to position 0), so temporarily
# extend the story past that point, plant the bookmark at the real target
# offset (now no longer "the very end"), and then trim the scratch text
# back off again - the bookmark stays put.
$d.Range($paraEnd, $paraEnd).InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $d.Range($paraEnd, $paraEnd)) | Out-Null
$d.Range($paraEnd, $paraEnd + 1).Delete()
